# Applies the "Actualización automática 2025-11-24 08:30:09" edits.
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# --- "VENTAS POR GRUPO" sheet ---
$wsGrupo.Range("M12").Value = 103.71
$wsGrupo.Range("L17").Value = 591.61
$wsGrupo.Range("L19").Value = "1 de 17"

# --- "VENTA MENSUAL" sheet ---
$wsMensual.Range("F12").Value = 103.71
$wsMensual.Range("F17").Value = 591.61
$wsMensual.Range("F19").Value = 2326.47
